# "Fruta / hortaliza, semanal" weekly update:
# Insert two new rows of data at the top of the data block (rows 5-6),
# pushing the existing rows 5-17 down to rows 7-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the current row 5 (shifts rows 5:17 down to 7:19,
# inheriting formatting/number-format from the row above, same as Excel UI).
$ws.Rows("5:6").Insert()

# New row 5
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 44881
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 300000000
$ws.Cells.Item(5, 7).Value = "Espárragos"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 2600
$ws.Cells.Item(5, 12).Value = 2700
$ws.Cells.Item(5, 13).Value = 2650
$ws.Cells.Item(5, 14).Value = "`$/kilo"
$ws.Cells.Item(5, 15).Value = "Provincia de Linares"
$ws.Cells.Item(5, 16).Value = 2650
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# New row 6
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value = "Bíobío"
$ws.Cells.Item(6, 4).Value = 44881
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 300000000
$ws.Cells.Item(6, 7).Value = "Espárragos"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Segunda"
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 2400
$ws.Cells.Item(6, 12).Value = 2400
$ws.Cells.Item(6, 13).Value = 2400
$ws.Cells.Item(6, 14).Value = "`$/kilo"
$ws.Cells.Item(6, 15).Value = "Provincia de Linares"
$ws.Cells.Item(6, 16).Value = 2400
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
